$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("meta")

# 1. Update the y_lab value (row 6, column B) to use double-semicolon separators.
$meta.Range("B6").Value = "-25;;-20;; -15;; -10;; -5;; 0;; 5"

# 2. Update the y_r_lab value (row 19, column B) to use double-semicolon separators.
$meta.Range("B19").Value = "-5;;-4;; -3;; -2;; -1;; 0;; 1"

# 3. Insert a new "style" / "default" row before the existing trailing blank row (row 23),
#    pushing the blank row down to row 24.
$meta.Rows("23:23").Insert()
$meta.Range("A23").Value = "style"
$meta.Range("B23").Value = "default"
